# Applies the CDS Input file update for TC02_CDS_Filter_FileType-BAM.xlsx:
# Replaces the "startup" sheet's Participants-tab Neo4j query (B2) with the
# revised Cypher query (adds diagnosis/genomic_info optional matches, a
# second pass to collect+sort sample ids, and sorts the final result).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['BAM']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$ws.Range("B2").Value = $newParticipantQuery

# The cell is wrap-text formatted, so the row grows from 10 wrapped lines to
# 18 -- mirror Excel's auto row-height recalculation for the new text.
$ws.Rows.Item(2).RowHeight = 279

# Match the author's final cursor position recorded in the saved file.
$ws.Range("C4").Select() | Out-Null
